$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look numeric keep their exact original
# text representation (e.g. trailing zeros) instead of being silently
# converted to Excel numbers, which would drop formatting.
$ws.Range('D2').Value = '65.579.64'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '2.660.64'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '600.64'
$ws.Range('E5').Value = '  -1.35%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '157.12'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.615'
$ws.Range('E8').Value = '  +4.43%  '
$ws.Range('E9').Value = '  -0.22%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '5.90'
$ws.Range('E10').Value = '  -0.33%  '
$ws.Range('E11').Value = '  -0.93%  '
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '29.33'
$ws.Range('E13').Value = '  -2.54%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000197'
$ws.Range('E14').Value = '  -1.48%  '
$ws.Range('D15').Value = '3.141.13'
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('D16').Value = '65.359.78'
$ws.Range('E16').Value = '  -0.30%  '
$ws.Range('D17').Value = '2.667.30'
$ws.Range('E17').Value = '  -0.48%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '12.66'
$ws.Range('E18').Value = '  -0.64%  '
$ws.Range('E19').Value = '  -2.37%  '
$ws.Range('E20').Value = '  +1.29%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '351.20'
$ws.Range('E21').Value = '  -2.58%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '69.41'
$ws.Range('E23').Value = '  -1.14%  '
$ws.Range('E24').Value = '  +4.40%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.76'
$ws.Range('E25').Value = '  +1.75%  '
$ws.Range('E26').Value = '  -4.99%  '
$ws.Range('E27').Value = '  -0.53%  '
$ws.Range('E28').Value = '  -3.23%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.09'
$ws.Range('E29').Value = '  -1.01%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '544.63'
$ws.Range('E30').Value = '  +2.54%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.997'
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.17'
$ws.Range('E32').Value = '  -2.28%  '
$ws.Range('E33').Value = '  -3.21%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.52'
$ws.Range('E34').Value = '  +1.61%  '
$ws.Range('E35').Value = '  -3.02%  '
$ws.Range('E36').Value = '  -2.83%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '20.38'
$ws.Range('E37').Value = '  -1.64%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '159.48'
$ws.Range('E38').Value = '  -2.22%  '
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.94'
$ws.Range('E40').Value = '  -2.72%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '42.66'
$ws.Range('E42').Value = '  +1.35%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '164.73'
$ws.Range('E43').Value = '  -0.71%  '
$ws.Range('E44').Value = '  -2.05%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.31'
$ws.Range('E45').Value = '  -1.29%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0611'
$ws.Range('E46').Value = '  -0.36%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '23.08'
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0259'
$ws.Range('E48').Value = '  -2.14%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.642'
$ws.Range('E49').Value = '  -2.17%  '
$ws.Range('E50').Value = '  +2.73%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '20.18'
$ws.Range('E51').Value = '  +1.43%  '
